$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert the new "2022-Q4" sheet right after "总计".
#    The current "2022-Q3" tab already holds the fund-level detail
#    rows that belong to the new quarter, so duplicate it (keeps all
#    formatting/styles byte-identical) and rename the copy.
# ------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)
$wsQ3 = $wb.Worksheets.Item("2022-Q3")
$wsQ3.Copy($null, $wsTotal)

$wsNew = $wb.Worksheets.Item(2)
$wsNew.Name = "2022-Q4"

# ------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: push every existing quarter
#    row down by one and fill in the brand-new 2022-Q4 figures at
#    the top; the row that falls off the bottom (2020-Q4) becomes a
#    new row 10.
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("总计")

for ($r = 9; $r -ge 2; $r--) {
    $ws1.Range("B$($r + 1)").Value = $ws1.Range("B$r").Value()
    $ws1.Range("C$($r + 1)").Value = $ws1.Range("C$r").Value()
    $ws1.Range("D$($r + 1)").Value = $ws1.Range("D$r").Value()
}

# Row 10 is brand new - copy row 9's formatting (thin border / bold /
# centered index-column look) onto it before writing the value.
$ws1.Range("A9").Copy()
$ws1.Range("A10").PasteSpecial(-4122)

$ws1.Range("B2").Value = "2022-Q4"
$ws1.Range("C2").Value = 5
$ws1.Range("D2").Value = 2.15

# Column A is simply the running 0-based row index.
for ($r = 2; $r -le 10; $r++) {
    $ws1.Range("A$r").Value = $r - 2
}
